# Commit #5: insurance, claim, debt, investment done
# - Rename sheet "債務" (debt) to "事業投資" (business investment)
# - Sheet "保險" (insurance): fix header row, replace policy-number column
#   with a constant "insurance" category column, and append the standard
#   property_category/category/date/legislator_name/legislator_id/
#   source_file/index metadata columns (F:K).
# - Sheet "事業投資" (ex "債務"): fix header row (English keys), drop the
#   buggy duplicated-header data row, and append the same metadata columns
#   (H:N) used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "債務" -> "事業投資"
# ---------------------------------------------------------------------
$wsInvest = $wb.Worksheets.Item("債務")
$wsInvest.Name = "事業投資"

# ---------------------------------------------------------------------
# 2. Sheet "保險" (insurance)
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Header row (row 1): B..K
$insHeaders = @("company","name","owner","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $insHeaders.Length; $i++) {
    $wsIns.Cells.Item(1, 2 + $i).Value = $insHeaders[$i]
}

# Data rows 2..28: column E becomes the constant "insurance" category,
# columns F..K get the standard metadata. Columns B (company), C (name)
# and D (owner) already hold the correct values.
for ($r = 2; $r -le 28; $r++) {
    $indexVal = $wsIns.Cells.Item($r, 1).Value2

    $wsIns.Cells.Item($r, 5).Value  = "insurance"          # E: property_category
    $wsIns.Cells.Item($r, 6).Value  = "normal"             # F: category
    $wsIns.Cells.Item($r, 7).Value  = "'2012-04-20"        # G: date (keep as text)
    $wsIns.Cells.Item($r, 8).Value  = "尤美女"              # H: legislator_name
    $wsIns.Cells.Item($r, 9).Value  = 1730                 # I: legislator_id
    $wsIns.Cells.Item($r, 10).Value = "tmp36451"           # J: source_file
    $wsIns.Cells.Item($r, 11).Value = $indexVal            # K: index
}

# ---------------------------------------------------------------------
# 3. Sheet "事業投資" (ex "債務" / investment)
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("事業投資")

# Row 2 is a buggy duplicate of the header row (A2=211) - remove it so the
# real data (previously rows 3 & 4, A=212/213) moves up to rows 2 & 3.
$wsInv.Rows.Item(2).Delete()

# Header row (row 1): B..G get English keys, H..N are new metadata columns
$invHeaders = @("owner","company","address","total","register_date","register_reason","property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $invHeaders.Length; $i++) {
    $wsInv.Cells.Item(1, 2 + $i).Value = $invHeaders[$i]
}

# Data rows 2..3: columns B..G already hold correct values (owner, company,
# address, total, register_date, register_reason); append H..N metadata.
for ($r = 2; $r -le 3; $r++) {
    $indexVal = $wsInv.Cells.Item($r, 1).Value2

    $wsInv.Cells.Item($r, 8).Value  = "investment"     # H: property_category
    $wsInv.Cells.Item($r, 9).Value  = "normal"         # I: category
    $wsInv.Cells.Item($r, 10).Value = "'2012-04-20"    # J: date (keep as text)
    $wsInv.Cells.Item($r, 11).Value = "尤美女"          # K: legislator_name
    $wsInv.Cells.Item($r, 12).Value = 1730             # L: legislator_id
    $wsInv.Cells.Item($r, 13).Value = "tmp36451"       # M: source_file
    $wsInv.Cells.Item($r, 14).Value = $indexVal        # N: index
}
